$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "1.00", "586.72").
# Force the cell format to Text first so Excel keeps the exact literal
# string instead of silently coercing it to a Double (dropping trailing
# zeros / introducing floating point noise).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.735.75"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.494.57"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.72"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.09"
$ws.Range("E6").Value = "  +4.26%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("E9").Value = "  +5.01%  "
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.340"
$ws.Range("E11").Value = "  +3.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.94"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.71"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.922.32"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.733.98"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.499.53"
$ws.Range("E17").Value = "  +3.84%  "
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.45"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "352.15"
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.73"
$ws.Range("E23").Value = "  +3.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.27"
$ws.Range("E24").Value = "  +2.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.78"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.20"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.622.00"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0912"
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "512.98"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.86"
$ws.Range("E31").Value = "  +3.94%  "
$ws.Range("E32").Value = "  +3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.79"
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.121"
$ws.Range("E35").Value = "  +7.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.00"
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.45"
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.70"
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  +5.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.331"
$ws.Range("E42").Value = "  +2.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.87"
$ws.Range("E43").Value = "  +2.87%  "
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "144.35"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("E46").Value = "  +2.88%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.515"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0258"
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("E50").Value = "  +2.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.586"
$ws.Range("E51").Value = "  +1.41%  "
